$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 25000
$ws.Range("F2").Value = 33000
$ws.Range("G2").Value = 100
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 31000
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 28000

$ws.Range("L2").Select()
